# Generate Report for Handback
#
# The workbook tracks localization handoff/handback status for each
# target language (zh-cn, de-de). This script records that the latest
# handoff package has been handed back and is in sync with en-US:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" columns are filled
#     in (re-using the same source .md / translated .xlf that were used
#     for the handoff), each as a hyperlink matching the existing ones
#   - The "Latest Handback DateTime" column is stamped with the handback
#     time

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- zh-cn sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column (B2): handoff -> handed back
$ws.Range("B2").Value = $statusText

# Latest Target File (E2) - same source markdown file as the handoff
$targetFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/308dc781a35ea0a35e2b66936506b528d0e623d4/e2e/a52e999f-6fe3-49f2-bb69-ed923be8f7d7.md"
$ws.Hyperlinks.Add($ws.Range("E2"), $targetFileUrl, "", "", "a52e999f-6fe3-49f2-bb69-ed923be8f7d7.md") | Out-Null

# Latest Handback File (F2) - the translated xlf that was handed back
$handbackFileUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/06c2079724bcf65805aa3272da1b70dff1e4c0c2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a52e999f-6fe3-49f2-bb69-ed923be8f7d7.878c4d101c65d61400847fbd7a13cd051d60a47f.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), $handbackFileUrl, "", "", "a52e999f-6fe3-49f2-bb69-ed923be8f7d7.878c4d101c65d61400847fbd7a13cd051d60a47f.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (G2)
$ws.Range("G2").Value = "2016-03-02 11:53:16"

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Status column (B2): handoff -> handed back
$ws.Range("B2").Value = $statusText

# Latest Target File (E2) - same source markdown file as the handoff
$targetFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/308dc781a35ea0a35e2b66936506b528d0e623d4/e2e/a52e999f-6fe3-49f2-bb69-ed923be8f7d7.md"
$ws.Hyperlinks.Add($ws.Range("E2"), $targetFileUrl, "", "", "a52e999f-6fe3-49f2-bb69-ed923be8f7d7.md") | Out-Null

# Latest Handback File (F2) - the translated xlf that was handed back
$handbackFileUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87194b21fb43f0a2e3f741d378d4c7047347039f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a52e999f-6fe3-49f2-bb69-ed923be8f7d7.878c4d101c65d61400847fbd7a13cd051d60a47f.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), $handbackFileUrl, "", "", "a52e999f-6fe3-49f2-bb69-ed923be8f7d7.878c4d101c65d61400847fbd7a13cd051d60a47f.de-de.xlf") | Out-Null

# Latest Handback DateTime (G2)
$ws.Range("G2").Value = "2016-03-02 11:53:36"
